$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "89.731.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.21%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.076.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.44%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +9.21%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "617.91"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.57%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -9.98%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.362"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.97%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.076.05"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.711"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.85%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.198"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.56%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.52%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.28"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.15%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.722.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.07%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.37"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.38%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.669.81"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.85%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.092.89"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.30%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.82"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.98%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000212"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.23%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.79"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.57%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "432.84"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -9.09%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.42"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.77"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.34%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.75"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.70%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "81.36"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -15.26%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.273.45"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.69%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.03%  "

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.08"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.99%  "

# Row 31
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.33%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.192"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.28%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.61"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.24%  "

# Row 35
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.152"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.69"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.46%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.09"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.86%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "495.19"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.43%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.16%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.26"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.44%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +55.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0879"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.10"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.58%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.396"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -5.86%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.71"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.76%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.86"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.17%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.675"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.33%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.34"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.63%  "

# Row 50
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.30"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.57%  "

# Row 51
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.22%  "
